# Add the "2022-Q3" quarter: a new per-fund detail sheet plus a new
# summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating the "2022-Q2"
#    sheet (so it inherits the same sheetPr/margins/header styling),
#    placing the copy immediately before it -> ends up in slot 2.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# The duplicated sheet has two fund rows; 2022-Q3 only has one, so drop
# row 3.
$q3Sheet.Rows.Item(3).Delete()

# Overwrite row 2 with the new quarter's single holding. Columns B:G
# hold text (fund code / name / numbers-as-text), column A/H are numeric.
$q3Sheet.Range("B2:G2").NumberFormat = "@"
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "540002"
$q3Sheet.Range("C2").Value = "汇丰晋信龙腾混合"
$q3Sheet.Range("D2").Value = "4.72"
$q3Sheet.Range("E2").Value = "93.98"
$q3Sheet.Range("F2").Value = "6.04"
$q3Sheet.Range("G2").Value = "0.2851"
$q3Sheet.Range("H2").Value = 6

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data (row 2) and push the existing quarters down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$rows = @(
    @(0, "2022-Q3", 1, 0.29),
    @(1, "2022-Q2", 2, 0.37),
    @(2, "2022-Q1", 2, 0.6),
    @(3, "2021-Q3", 1, 0.02),
    @(4, "2021-Q2", 3, 0.03),
    @(5, "2021-Q1", 3, 0.06)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $vals = $rows[$i]
    $total.Cells.Item($r, 1).Value = $vals[0]
    $total.Cells.Item($r, 2).Value = $vals[1]
    $total.Cells.Item($r, 3).Value = $vals[2]
    $total.Cells.Item($r, 4).Value = $vals[3]
}

# Row 7 is brand new - copy the formatting (bold/border style) of column
# A from the row above so it matches the rest of the index column.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Restore the original active tab (the last sheet, "2021-Q1").
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
